$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 523; this shifts existing rows 523:632 down to 524:633
$ws.Rows.Item(523).Insert()

# Populate the newly inserted row 523 with its data
$ws.Cells.Item(523, 1).Value = 3
$ws.Cells.Item(523, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(523, 3).Value = "Coquimbo"
$ws.Cells.Item(523, 4).Value = 45211
$ws.Cells.Item(523, 5).Value = 5
$ws.Cells.Item(523, 6).Value = 100112009
$ws.Cells.Item(523, 7).Value = "Acelga"
$ws.Cells.Item(523, 8).Value = "Sin especificar"
$ws.Cells.Item(523, 9).Value = "Primera"
$ws.Cells.Item(523, 10).Value = 240
$ws.Cells.Item(523, 11).Value = 3000
$ws.Cells.Item(523, 12).Value = 3500
$ws.Cells.Item(523, 13).Value = 3250
$ws.Cells.Item(523, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(523, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(523, 16).Value = 542
$ws.Cells.Item(523, 17).Value = 6
$ws.Cells.Item(523, 18).Value = "Hortaliza"
